$d = $word.ActiveDocument

# 1) Insert a leading space before "Lab: " at the very start of the document.
$start = $d.Range(0, 0)
$start.InsertBefore(" ")

# 2) Merge the hyperlink runs "h" + "e" + "re" into a single run "here"
#    by doing a find & replace on the unique text "here" within the hyperlink.
$rng = $d.Content
$rng.Find.Execute("here", $false, $false, $false, $false, $false, $true, 1, $false, "here", 2)

Write-Output "done"
